$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.965.88"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").Value = "3.171.18"
$ws.Range("E3").Value = "  +4.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'579.70"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").Value = "'151.49"
$ws.Range("E6").Value = "  +7.18%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.170.86"
$ws.Range("E8").Value = "  +4.31%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("E10").Value = "  +6.61%  "
$ws.Range("D11").Value = "'6.21"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "'0.503"
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +18.58%  "
$ws.Range("D14").Value = "'37.58"
$ws.Range("E14").Value = "  +6.20%  "
$ws.Range("D15").Value = "3.692.87"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").Value = "65.060.93"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.172.19"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.18"
$ws.Range("E18").Value = "  +6.47%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'511.49"
$ws.Range("E20").Value = "  +6.55%  "
$ws.Range("D21").Value = "'14.88"
$ws.Range("E21").Value = "  +5.91%  "
$ws.Range("E22").Value = "  +6.93%  "
$ws.Range("D23").Value = "'15.53"
$ws.Range("E23").Value = "  +6.68%  "
$ws.Range("D24").Value = "'7.82"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("D25").Value = "'85.15"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'9.13"
$ws.Range("E27").Value = "  +13.47%  "
$ws.Range("E28").Value = "  +5.22%  "
$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  +8.71%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'27.83"
$ws.Range("E30").Value = "  +6.80%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.80"
$ws.Range("E31").Value = "  +15.35%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("E34").Value = "  +12.54%  "
$ws.Range("E35").Value = "  +6.57%  "
$ws.Range("D36").Value = "'55.71"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "'0.0905"
$ws.Range("E37").Value = "  +11.65%  "
$ws.Range("D38").Value = "'475.20"
$ws.Range("E38").Value = "  +7.93%  "
$ws.Range("D39").Value = "'3.10"
$ws.Range("E39").Value = "  +12.87%  "
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("D41").Value = "'8.67"
$ws.Range("E41").Value = "  +4.68%  "
$ws.Range("D42").Value = "3.074.12"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.45"
$ws.Range("E44").Value = "  +10.14%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.285"
$ws.Range("E45").Value = "  +6.15%  "
$ws.Range("D46").Value = "'29.46"
$ws.Range("E46").Value = "  +6.62%  "
$ws.Range("D47").Value = "0.0₃0608"
$ws.Range("E47").Value = "  +19.45%  "
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("D50").Value = "'2.26"
$ws.Range("E50").Value = "  +8.92%  "
$ws.Range("D51").Value = "'120.49"
$ws.Range("E51").Value = "  +2.40%  "
